# Coastal Surface Piercing Profilers - Updated Coastal CSPP ingest and cal sheets
# - corrected instrument reference designators (GP05MOAS-GL003 -> GP05MOAS-GL453)
# - Asset_Cal_Info deployment numbers corrected from 3 to 1
# - active sheet/selection moved to Asset_Cal_Info

$wb = $excel.ActiveWorkbook
$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsCal = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Moorings sheet: fix reference designator and deployment number ---
$wsMoorings.Range("A2").Value = "GP05MOAS-GL453"
$wsMoorings.Range("C2").Value = 1

# --- Asset_Cal_Info sheet: fix reference designators and deployment numbers ---
$wsCal.Range("A3").Value = "GP05MOAS-GL453-00-ENG000000"
$wsCal.Range("A4").Value = "GP05MOAS-GL453-01-FLORDM000"
$wsCal.Range("A5").Value = "GP05MOAS-GL453-01-FLORDM000"
$wsCal.Range("A6").Value = "GP05MOAS-GL453-01-FLORDM000"
$wsCal.Range("A7").Value = "GP05MOAS-GL453-01-FLORDM000"
$wsCal.Range("A8").Value = "GP05MOAS-GL453-02-DOSTAM000"
$wsCal.Range("A9").Value = "GP05MOAS-GL453-04-CTDGVM000"

$wsCal.Range("C3:C9").Value = 1

# --- Update view state: Moorings selection moves, Asset_Cal_Info becomes active tab ---
$wsMoorings.Range("E22").Select()
$wsCal.Activate()
$wsCal.Range("C10").Select()
